$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = "Â±"
$goodChar = "±"

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val.GetType().Name -eq "String" -and $val.Contains($badChar)) {
            $cell.Value2 = $val.Replace($badChar, $goodChar)
        }
    }
}
